$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New hiking-journal entries (visited by me / image flags sync).
# Values are written first (so the E1/G1 summary formulas correctly
# pick the new cells up as precedents), and the row formatting
# (matching the previous group, rows 54-57) is copied in afterwards.
# ---------------------------------------------------------------------

# --- Group 1: rows 58-61 (trip on 2025-03-08) ---
$ws.Range("A58").Value = 99
$ws.Range("B58").Value = "HWN099 Harzer Holzdampflok"
$ws.Range("C58").Value = "03/08/2025"
$ws.Range("D58").Value = "https://www.komoot.com/de-de/tour/2088114497"
$ws.Range("E58").Value = 1

$ws.Range("A59").Value = 100
$ws.Range("B59").Value = "HWN100 Ebersburg"
$ws.Range("C59").Value = "03/08/2025"

$ws.Range("A60").Value = 218
$ws.Range("B60").Value = "HWN218 Neustädter Talsperre"
$ws.Range("C60").Value = "03/08/2025"

$ws.Range("A61").Value = 98
$ws.Range("B61").Value = "HWN098 Ruine Hohnstein"
$ws.Range("C61").Value = "03/08/2025"

$ws.Range("A54:E57").Copy()
$ws.Range("A58").PasteSpecial(-4122)

$ws.Range("D58:D61").Merge()
$ws.Range("E58:E61").Merge()
$ws.Hyperlinks.Add($ws.Range("D58"), "https://www.komoot.com/de-de/tour/2088114497")
# Adding the hyperlink overwrites the cell format with Excel's built-in
# "Hyperlink" style; restore the journal's own link style (same as D54).
$ws.Range("D54").Copy()
$ws.Range("D58").PasteSpecial(-4122)

# --- Group 2: rows 62-64 (trip on 2025-03-23) ---
$ws.Range("A62").Value = 164
$ws.Range("B62").Value = "HWN164 Stiefmutter"
$ws.Range("C62").Value = "03/23/2025"
$ws.Range("D62").Value = "https://www.komoot.com/de-de/tour/2112900887"
$ws.Range("E62").Value = 1

$ws.Range("A63").Value = 165
$ws.Range("B63").Value = "HWN165 Wendel-Eiche"
$ws.Range("C63").Value = "03/23/2025"

$ws.Range("A64").Value = 90
$ws.Range("B64").Value = "HWN090 Roter Schuss"
$ws.Range("C64").Value = "03/23/2025"

$ws.Range("A54:E56").Copy()
$ws.Range("A62").PasteSpecial(-4122)

$ws.Range("D62:D64").Merge()
$ws.Range("E62:E64").Merge()
$ws.Hyperlinks.Add($ws.Range("D62"), "https://www.komoot.com/de-de/tour/2112900887")
$ws.Range("D54").Copy()
$ws.Range("D62").PasteSpecial(-4122)

# --- Group 3: rows 65-68 (trip on 2025-05-01) ---
$ws.Range("A65").Value = 35
$ws.Range("B65").Value = "HWN035 Gasthaus Armeleuteberg"
$ws.Range("C65").Value = "05/01/2025"
$ws.Range("D65").Value = "https://www.komoot.com/de-de/tour/2021000530"
$ws.Range("E65").Value = 1

$ws.Range("A66").Value = 34
$ws.Range("B66").Value = "HWN034 Scharfenstein"
$ws.Range("C66").Value = "05/01/2025"

$ws.Range("A67").Value = 32
$ws.Range("B67").Value = "HWN032 Gasthaus Christianental"
$ws.Range("C67").Value = "05/01/2025"

$ws.Range("A68").Value = 31
$ws.Range("B68").Value = "HWN031 Agnesberg"
$ws.Range("C68").Value = "05/01/2025"

$ws.Range("A54:E57").Copy()
$ws.Range("A65").PasteSpecial(-4122)

$ws.Range("D65:D68").Merge()
$ws.Range("E65:E68").Merge()
$ws.Hyperlinks.Add($ws.Range("D65"), "https://www.komoot.com/de-de/tour/2021000530")
$ws.Range("D54").Copy()
$ws.Range("D65").PasteSpecial(-4122)

# Recalculate SUM(E2:E126) and COUNTA(A2:A250) formulas in row 1.
$excel.Calculate()
